# Update regression coefficients on HU_U1a (the "U1a" matrix referenced in the
# commit message) and flip which sheet/tab is active: HU_U1a becomes selected
# instead of HU_U2b ("HO1a"). The B2:N13 data block also loses its explicit
# number-format styling (reverts to the default "Normal" style).
$wb  = $excel.ActiveWorkbook
$ws1a = $wb.Worksheets.Item("HU_U1a")
$ws2b = $wb.Worksheets.Item("HU_U2b")

# The data block B2:N13 drops its custom cell style (numFmt + font) and reverts
# to the workbook default "Normal" style.
$ws1a.Range("B2:N13").Style = "Normal"

# Row 2
$ws1a.Range("B2").Value = -0.14207781397738986
$ws1a.Range("C2").Value = 0.039099736512404268
$ws1a.Range("D2").Value = -0.02618773595208243
$ws1a.Range("E2").Value = 0.00058598518496328841
$ws1a.Range("F2").Value = 0.0037526975228583531
$ws1a.Range("G2").Value = 0.010572877447602047
$ws1a.Range("H2").Value = 0.023104007605285824
$ws1a.Range("I2").Value = 0.003123338556850451
$ws1a.Range("J2").Value = 0.00040560258119748914
$ws1a.Range("K2").Value = -0.041727704056742893
$ws1a.Range("L2").Value = -0.033616158247294754
$ws1a.Range("M2").Value = 0.0014633370720082244
$ws1a.Range("N2").Value = 0.26993629683316778

# Row 3
$ws1a.Range("B3").Value = 0.91595425602714586
$ws1a.Range("C3").Value = -0.02618773595208243
$ws1a.Range("D3").Value = 0.77227331779013397
$ws1a.Range("E3").Value = -0.017868253106205989
$ws1a.Range("F3").Value = -0.037394954911670553
$ws1a.Range("G3").Value = -0.13967597249802388
$ws1a.Range("H3").Value = 0.066725660850274315
$ws1a.Range("I3").Value = -0.0080128946048574345
$ws1a.Range("J3").Value = -0.047346964551454485
$ws1a.Range("K3").Value = -0.10753579576157302
$ws1a.Range("L3").Value = -0.016805353605709783
$ws1a.Range("M3").Value = 0.0014029596075838535
$ws1a.Range("N3").Value = -7.9691541135136958

# Row 4
$ws1a.Range("B4").Value = -0.022270470896169778
$ws1a.Range("C4").Value = 0.00058598518496328841
$ws1a.Range("D4").Value = -0.017868253106205989
$ws1a.Range("E4").Value = 0.00041502823114809503
$ws1a.Range("F4").Value = 0.00072805086520663642
$ws1a.Range("G4").Value = 0.00308056731284026
$ws1a.Range("H4").Value = -0.0015257164527659236
$ws1a.Range("I4").Value = 0.00025399765414441458
$ws1a.Range("J4").Value = 0.0011689079215027876
$ws1a.Range("K4").Value = 0.0026427450608180192
$ws1a.Range("L4").Value = 0.00033224448668943418
$ws1a.Range("M4").Value = -0.0000091271085454844997
$ws1a.Range("N4").Value = 0.18306269763041705

# Row 5
$ws1a.Range("B5").Value = -0.32709430964378433
$ws1a.Range("C5").Value = 0.0037526975228583531
$ws1a.Range("D5").Value = -0.037394954911670553
$ws1a.Range("E5").Value = 0.00072805086520663642
$ws1a.Range("F5").Value = 0.070513101545748871
$ws1a.Range("G5").Value = 0.01996431402641629
$ws1a.Range("H5").Value = 0.0061561442211249453
$ws1a.Range("I5").Value = -0.0033204129319544164
$ws1a.Range("J5").Value = -0.0053198194806704071
$ws1a.Range("K5").Value = -0.003791241687210989
$ws1a.Range("L5").Value = 0.0064611284148303316
$ws1a.Range("M5").Value = 0.0039370561063523622
$ws1a.Range("N5").Value = 0.39982721517354414

# Row 6
$ws1a.Range("B6").Value = -0.15178401145857037
$ws1a.Range("C6").Value = 0.010572877447602047
$ws1a.Range("D6").Value = -0.13967597249802388
$ws1a.Range("E6").Value = 0.00308056731284026
$ws1a.Range("F6").Value = 0.01996431402641629
$ws1a.Range("G6").Value = 0.041233758580018598
$ws1a.Range("H6").Value = -0.01822709277687808
$ws1a.Range("I6").Value = -0.0044678536908201588
$ws1a.Range("J6").Value = 0.0040601149925975979
$ws1a.Range("K6").Value = -0.0026529309933064904
$ws1a.Range("L6").Value = 0.0014722325361961608
$ws1a.Range("M6").Value = -0.0021085473478027611
$ws1a.Range("N6").Value = 1.5538975083032049

# Row 7
$ws1a.Range("B7").Value = -0.062890788086950353
$ws1a.Range("C7").Value = 0.023104007605285824
$ws1a.Range("D7").Value = 0.066725660850274315
$ws1a.Range("E7").Value = -0.0015257164527659236
$ws1a.Range("F7").Value = 0.0061561442211249453
$ws1a.Range("G7").Value = -0.01822709277687808
$ws1a.Range("H7").Value = 0.10911202882183177
$ws1a.Range("I7").Value = -0.0054740086371200536
$ws1a.Range("J7").Value = -0.032259892463744194
$ws1a.Range("K7").Value = -0.0070095466919185556
$ws1a.Range("L7").Value = 0.0071115188747810087
$ws1a.Range("M7").Value = 0.0013879288912782145
$ws1a.Range("N7").Value = -0.61313481186356

# Row 8
$ws1a.Range("B8").Value = 0.080314423940431381
$ws1a.Range("C8").Value = 0.003123338556850451
$ws1a.Range("D8").Value = -0.0080128946048574345
$ws1a.Range("E8").Value = 0.00025399765414441458
$ws1a.Range("F8").Value = -0.0033204129319544164
$ws1a.Range("G8").Value = -0.0044678536908201588
$ws1a.Range("H8").Value = -0.0054740086371200536
$ws1a.Range("I8").Value = 0.0062801808381461122
$ws1a.Range("J8").Value = 0.0064121021336828142
$ws1a.Range("K8").Value = 0.004028493455235406
$ws1a.Range("L8").Value = -0.0043032014843607316
$ws1a.Range("M8").Value = 0.001159068718196727
$ws1a.Range("N8").Value = 0.0083563390140843907

# Row 9
$ws1a.Range("B9").Value = -0.069952798698799626
$ws1a.Range("C9").Value = 0.00040560258119748914
$ws1a.Range("D9").Value = -0.047346964551454485
$ws1a.Range("E9").Value = 0.0011689079215027876
$ws1a.Range("F9").Value = -0.0053198194806704071
$ws1a.Range("G9").Value = 0.0040601149925975979
$ws1a.Range("H9").Value = -0.032259892463744194
$ws1a.Range("I9").Value = 0.0064121021336828142
$ws1a.Range("J9").Value = 0.016737684879381386
$ws1a.Range("K9").Value = 0.0007795262747052889
$ws1a.Range("L9").Value = -0.016463150465742392
$ws1a.Range("M9").Value = 0.0015430686955089073
$ws1a.Range("N9").Value = 0.38657918714809858

# Row 10
$ws1a.Range("B10").Value = 0.28698026707993862
$ws1a.Range("C10").Value = -0.041727704056742893
$ws1a.Range("D10").Value = -0.10753579576157302
$ws1a.Range("E10").Value = 0.0026427450608180192
$ws1a.Range("F10").Value = -0.003791241687210989
$ws1a.Range("G10").Value = -0.0026529309933064904
$ws1a.Range("H10").Value = -0.0070095466919185556
$ws1a.Range("I10").Value = 0.004028493455235406
$ws1a.Range("J10").Value = 0.0007795262747052889
$ws1a.Range("K10").Value = 0.11758966698153482
$ws1a.Range("L10").Value = 0.076596526702307391
$ws1a.Range("M10").Value = -0.0023975540522050975
$ws1a.Range("N10").Value = 1.0528365402341322

# Row 11
$ws1a.Range("B11").Value = 0.4029488529453894
$ws1a.Range("C11").Value = -0.033616158247294754
$ws1a.Range("D11").Value = -0.016805353605709783
$ws1a.Range("E11").Value = 0.00033224448668943418
$ws1a.Range("F11").Value = 0.0064611284148303316
$ws1a.Range("G11").Value = 0.0014722325361961608
$ws1a.Range("H11").Value = 0.0071115188747810087
$ws1a.Range("I11").Value = -0.0043032014843607316
$ws1a.Range("J11").Value = -0.016463150465742392
$ws1a.Range("K11").Value = 0.076596526702307391
$ws1a.Range("L11").Value = 0.090075935470496415
$ws1a.Range("M11").Value = -0.0056347966817148973
$ws1a.Range("N11").Value = 0.29678918708959046

# Row 12
$ws1a.Range("B12").Value = -0.06443959417262439
$ws1a.Range("C12").Value = 0.0014633370720082244
$ws1a.Range("D12").Value = 0.0014029596075838535
$ws1a.Range("E12").Value = -0.0000091271085454844997
$ws1a.Range("F12").Value = 0.0039370561063523622
$ws1a.Range("G12").Value = -0.0021085473478027611
$ws1a.Range("H12").Value = 0.0013879288912782145
$ws1a.Range("I12").Value = 0.001159068718196727
$ws1a.Range("J12").Value = 0.0015430686955089073
$ws1a.Range("K12").Value = -0.0023975540522050975
$ws1a.Range("L12").Value = -0.0056347966817148973
$ws1a.Range("M12").Value = 0.00145632446695021
$ws1a.Range("N12").Value = -0.050505140231000029

# Row 13
$ws1a.Range("B13").Value = -11.295675051520156
$ws1a.Range("C13").Value = 0.26993629683316778
$ws1a.Range("D13").Value = -7.9691541135136958
$ws1a.Range("E13").Value = 0.18306269763041705
$ws1a.Range("F13").Value = 0.39982721517354414
$ws1a.Range("G13").Value = 1.5538975083032049
$ws1a.Range("H13").Value = -0.61313481186356
$ws1a.Range("I13").Value = 0.0083563390140843907
$ws1a.Range("J13").Value = 0.38657918714809858
$ws1a.Range("K13").Value = 1.0528365402341322
$ws1a.Range("L13").Value = 0.29678918708959046
$ws1a.Range("M13").Value = -0.050505140231000029
$ws1a.Range("N13").Value = 83.733272096101246

# HU_U1a becomes the active/selected tab (was HU_U2b before the edit).
$ws1a.Activate()
